$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 11
$ws.Range("E3").Value = 11
$ws.Range("E4").Value = 11
$ws.Range("E5").Value = 7
$ws.Range("E6").Value = 11
$ws.Range("E7").Value = 7
$ws.Range("E8").Value = 11
$ws.Range("E9").Value = 7
$ws.Range("E10").Value = 4
$ws.Range("E11").Value = 11
$ws.Range("E12").Value = 7
$ws.Range("E13").Value = 11
$ws.Range("E14").Value = 11
$ws.Range("E15").Value = 11
$ws.Range("E16").Value = 1
$ws.Range("E17").Value = 7
$ws.Range("E18").Value = 10
$ws.Range("F18").Value = 20251009
$ws.Range("E19").Value = 10
$ws.Range("F19").Value = 20251009
$ws.Range("E20").Value = 10
$ws.Range("F20").Value = 20251009
$ws.Range("E21").Value = 10
$ws.Range("F21").Value = 20251009
$ws.Range("E22").Value = 7
$ws.Range("E23").Value = 7
$ws.Range("E24").Value = 7
$ws.Range("E25").Value = 7
$ws.Range("E26").Value = 7
$ws.Range("E27").Value = 5
$ws.Range("E28").Value = 10
$ws.Range("F28").Value = 20251009
$ws.Range("E29").Value = 10
$ws.Range("F29").Value = 20251009
$ws.Range("E30").Value = 10
$ws.Range("F30").Value = 20251009
$ws.Range("E31").Value = 10
$ws.Range("F31").Value = 20251009
$ws.Range("E32").Value = 10
$ws.Range("F32").Value = 20251009
$ws.Range("E33").Value = 10
$ws.Range("F33").Value = 20251009
$ws.Range("E34").Value = 10
$ws.Range("F34").Value = 20251009
$ws.Range("E35").Value = 10
$ws.Range("F35").Value = 20251009
$ws.Range("E37").Value = 10
$ws.Range("F37").Value = 20251009
$ws.Range("E38").Value = 10
$ws.Range("F38").Value = 20251009
$ws.Range("E39").Value = 10
$ws.Range("F39").Value = 20251009
$ws.Range("E40").Value = 4
$ws.Range("E41").Value = 4
$ws.Range("E42").Value = 10
$ws.Range("F42").Value = 20251009
$ws.Range("E43").Value = 7
$ws.Range("E44").Value = 4
$ws.Range("E45").Value = 7
$ws.Range("E46").Value = 4
$ws.Range("E47").Value = 10
$ws.Range("F47").Value = 20251009
$ws.Range("E48").Value = 4
$ws.Range("E49").Value = 5
$ws.Range("E50").Value = 5
$ws.Range("E51").Value = 5
$ws.Range("E52").Value = 5
$ws.Range("E53").Value = 5
$ws.Range("E54").Value = 5
$ws.Range("E55").Value = 5
$ws.Range("E56").Value = 5
$ws.Range("E57").Value = 5
$ws.Range("E58").Value = 9
$ws.Range("E59").Value = 9
$ws.Range("E60").Value = 9
$ws.Range("E61").Value = 5
$ws.Range("E62").Value = 9
$ws.Range("E63").Value = 9
$ws.Range("E64").Value = 9
$ws.Range("E65").Value = 10
$ws.Range("F65").Value = 20251009
$ws.Range("E66").Value = 10
$ws.Range("F66").Value = 20251009
$ws.Range("E67").Value = 10
$ws.Range("F67").Value = 20251009
$ws.Range("E68").Value = 10
$ws.Range("F68").Value = 20251009
$ws.Range("E69").Value = 10
$ws.Range("F69").Value = 20251009
$ws.Range("E70").Value = 1
$ws.Range("E71").Value = 1
$ws.Range("E72").Value = 1
$ws.Range("E73").Value = 1
$ws.Range("E74").Value = 1
$ws.Range("E75").Value = 1
$ws.Range("E76").Value = 1
$ws.Range("E77").Value = 4
$ws.Range("E78").Value = 4
$ws.Range("E79").Value = 4
$ws.Range("E80").Value = 4
$ws.Range("E81").Value = 4
$ws.Range("E82").Value = 4
$ws.Range("E83").Value = 4
$ws.Range("E84").Value = 4
$ws.Range("E85").Value = 4
$ws.Range("E86").Value = 4
$ws.Range("E87").Value = 4
$ws.Range("E88").Value = 4
$ws.Range("E89").Value = 4
$ws.Range("E90").Value = 4
$ws.Range("E91").Value = 7
$ws.Range("E92").Value = 4
$ws.Range("E93").Value = 4
$ws.Range("E94").Value = 7
$ws.Range("F94").Value = 20251009
$ws.Range("E95").Value = 3
